$p = $ppt.ActivePresentation

# Insert a new slide right after the existing "Massive parallelization"
# title slide. Layout 2 (ppLayoutText) is PowerPoint's "Title and
# Content" layout -- a title placeholder plus a body/content
# placeholder -- matching the shapes on the new slide.
$s2 = $p.Slides.Add(2, 2)

# Give the new slide its title ("added slide titles for all slides").
$s2.Shapes.Title.TextFrame.TextRange.Text = "MapReduce"
